$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 7813.154
$ws.Range("I9").Value = 10120.6
$ws.Range("J9").Value = 121.666664
$ws.Range("K9").Value = 10120.6
$ws.Range("L9").Value = 121.666664
$ws.Range("M9").Value = -9951.6
$ws.Range("N9").Value = -459.666664
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H19").Value = 1159.6666
$ws.Range("I19").Value = 1194.5555
$ws.Range("J19").Value = 1124.7778
$ws.Range("K19").Value = 1194.5555
$ws.Range("L19").Value = 1124.7778
$ws.Range("M19").Value = -1019.5555
$ws.Range("N19").Value = -1474.7778
$ws.Range("H53").Value = 1320.5714
$ws.Range("I53").Value = 1067.5
$ws.Range("J53").Value = 1658
$ws.Range("K53").Value = 1067.5
$ws.Range("L53").Value = 1658
$ws.Range("M53").Value = -430.5
$ws.Range("N53").Value = -2932
$ws.Range("H70").Value = 3661.848
$ws.Range("I70").Value = 1383.6471
$ws.Range("J70").Value = 10116.75
$ws.Range("K70").Value = 4150.9413
$ws.Range("L70").Value = 30350.25
$ws.Range("M70").Value = -3880.9413
$ws.Range("N70").Value = -30890.25
$ws.Range("H73").Value = 3661.848
$ws.Range("I73").Value = 1383.6471
$ws.Range("J73").Value = 10116.75
$ws.Range("K73").Value = 4150.9413
$ws.Range("L73").Value = 30350.25
$ws.Range("M73").Value = -3214.9413
$ws.Range("N73").Value = -32222.25
$ws.Range("H103").Value = 811
$ws.Range("J103").Value = 775.25
$ws.Range("L103").Value = 2325.75
$ws.Range("N103").Value = -3497.75
$ws.Range("H132").Value = 2881.111
$ws.Range("I132").Value = 2725.9321
$ws.Range("K132").Value = 8177.7963
$ws.Range("M132").Value = -5647.7963
$ws.Range("H134").Value = 34998.734
$ws.Range("J134").Value = 34998.734
$ws.Range("L134").Value = 34998.734
$ws.Range("N134").Value = -45138.734

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1410.5555
$ws.Range("I2").Value = 972
$ws.Range("J2").Value = 1761.4
$ws.Range("K2").Value = 972
$ws.Range("L2").Value = 1761.4
$ws.Range("M2").Value = -859
$ws.Range("N2").Value = -1987.4
$ws.Range("H45").Value = 3679
$ws.Range("I45").Value = 3678.8
$ws.Range("K45").Value = 3678.8
$ws.Range("M45").Value = -3301.8
$ws.Range("H97").Value = 1064.909
$ws.Range("I97").Value = 580.6316
$ws.Range("K97").Value = 580.6316
$ws.Range("M97").Value = -84.63160000000005
$ws.Range("H116").Value = 1410.5555
$ws.Range("I116").Value = 972
$ws.Range("J116").Value = 1761.4
$ws.Range("K116").Value = 972
$ws.Range("L116").Value = 1761.4
$ws.Range("M116").Value = 1322
$ws.Range("N116").Value = -6349.4
$ws.Range("H122").Value = 2578.1
$ws.Range("I122").Value = 1412.15
$ws.Range("K122").Value = 4236.450000000001
$ws.Range("M122").Value = -1786.450000000001
$ws.Range("H132").Value = 2186.842
$ws.Range("I132").Value = 2014.7142
$ws.Range("J132").Value = 2668.8
$ws.Range("K132").Value = 6044.142599999999
$ws.Range("L132").Value = 8006.400000000001
$ws.Range("M132").Value = -3514.142599999999
$ws.Range("N132").Value = -13066.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1410.5555
$ws.Range("I3").Value = 972
$ws.Range("J3").Value = 1761.4
$ws.Range("K3").Value = 972
$ws.Range("L3").Value = 1761.4
$ws.Range("M3").Value = -858
$ws.Range("N3").Value = -1989.4
$ws.Range("H132").Value = 40000
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 86999
$ws.Range("J133").Value = 86999
$ws.Range("L133").Value = 86999
$ws.Range("N133").Value = -97119
$ws.Range("H134").Value = 4101.904
$ws.Range("I134").Value = 3540.5117
$ws.Range("J134").Value = 6784.1113
$ws.Range("K134").Value = 10621.5351
$ws.Range("L134").Value = 20352.3339
$ws.Range("M134").Value = -8086.535100000001
$ws.Range("N134").Value = -25422.3339
$ws.Range("H139").Value = 223249.75
$ws.Range("J139").Value = 223249.75
$ws.Range("L139").Value = 223249.75
$ws.Range("N139").Value = -233529.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 849.75
$ws.Range("I22").Value = 399.5
$ws.Range("J22").Value = 1300
$ws.Range("K22").Value = 399.5
$ws.Range("L22").Value = 1300
$ws.Range("M22").Value = -49.5
$ws.Range("N22").Value = -2000
$ws.Range("H62").Value = 6899.5
$ws.Range("I62").Value = 5300
$ws.Range("J62").Value = 8499
$ws.Range("K62").Value = 5300
$ws.Range("L62").Value = 8499
$ws.Range("M62").Value = -4676
$ws.Range("N62").Value = -9747
$ws.Range("H65").Value = 6899.5
$ws.Range("I65").Value = 5300
$ws.Range("J65").Value = 8499
$ws.Range("K65").Value = 26500
$ws.Range("L65").Value = 42495
$ws.Range("M65").Value = -23380
$ws.Range("N65").Value = -48735
$ws.Range("H105").Value = 1273.1
$ws.Range("I105").Value = 1273.1
$ws.Range("K105").Value = 1273.1
$ws.Range("M105").Value = 473.9000000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1703.9166
$ws.Range("I80").Value = 1550.4
$ws.Range("J80").Value = 1813.5714
$ws.Range("K80").Value = 1550.4
$ws.Range("L80").Value = 1813.5714
$ws.Range("M80").Value = -552.4000000000001
$ws.Range("N80").Value = -3809.5714
$ws.Range("H83").Value = 1703.9166
$ws.Range("I83").Value = 1550.4
$ws.Range("J83").Value = 1813.5714
$ws.Range("K83").Value = 7752
$ws.Range("L83").Value = 9067.857
$ws.Range("M83").Value = -2760
$ws.Range("N83").Value = -19051.857
$ws.Range("H135").Value = 94979.39999999999
$ws.Range("J135").Value = 94979.39999999999
$ws.Range("L135").Value = 94979.39999999999
$ws.Range("N135").Value = -105119.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H113").Value = 1532.8667
$ws.Range("I113").Value = 1282.4445
$ws.Range("K113").Value = 3847.3335
$ws.Range("M113").Value = -1677.3335
